$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.666448
$ws.Range("H2").Value = 3.332896
$ws.Range("I2").Value = 0.08698774157534103
$ws.Range("J2").Value = 0.0653244391585683
$ws.Range("M2").Value = 0.07111099999999999
$ws.Range("N2").Value = 0.142222
$ws.Range("O2").Value = 0.02711460746047303
$ws.Range("P2").Value = 0.02678527708115022
$ws.Range("Q2").Value = 0.118502783728
$ws.Range("R2").Value = 0.4740111349119999
$ws.Range("S2").Value = 0.002358638466688442
$ws.Range("T2").Value = 0.001749733203032992
$ws.Range("G3").Value = 1.666448
$ws.Range("H3").Value = 3.332896
$ws.Range("I3").Value = 0.08698774157534103
$ws.Range("J3").Value = 0.0653244391585683
$ws.Range("O3").Value = 0.0245904030281302
$ws.Range("P3").Value = 0.03643759694506741
$ws.Range("Q3").Value = 0.107470897968
$ws.Range("R3").Value = 0.644825387808
$ws.Range("S3").Value = 0.002139063623844474
$ws.Range("T3").Value = 0.00238026558472249
$ws.Range("G4").Value = 1.666448
$ws.Range("H4").Value = 3.332896
$ws.Range("I4").Value = 0.08698774157534103
$ws.Range("J4").Value = 0.0653244391585683
$ws.Range("M4").Value = 2.4870065
$ws.Range("N4").Value = 4.974013
$ws.Range("O4").Value = 0.9482949895113968
$ws.Range("P4").Value = 0.9367771259737823
$ws.Range("Q4").Value = 4.144467007912
$ws.Range("R4").Value = 16.577868031648
$ws.Range("S4").Value = 0.08249003948480813
$ws.Range("T4").Value = 0.06119444037081281
$ws.Range("I5").Value = 0.1669502665149541
$ws.Range("J5").Value = 0.1880598173367416
$ws.Range("M5").Value = 0.07111099999999999
$ws.Range("N5").Value = 0.142222
$ws.Range("O5").Value = 0.02711460746047303
$ws.Range("P5").Value = 0.02678527708115022
$ws.Range("Q5").Value = 0.227435164632
$ws.Range("R5").Value = 1.364610987792
$ws.Range("S5").Value = 0.004526790941974335
$ws.Range("T5").Value = 0.005037234315195122
$ws.Range("I6").Value = 0.1669502665149541
$ws.Range("J6").Value = 0.1880598173367416
$ws.Range("O6").Value = 0.0245904030281302
$ws.Range("P6").Value = 0.03643759694506741
$ws.Range("S6").Value = 0.00410537433925647
$ws.Range("T6").Value = 0.006852447825679192
$ws.Range("I7").Value = 0.1669502665149541
$ws.Range("J7").Value = 0.1880598173367416
$ws.Range("M7").Value = 2.4870065
$ws.Range("N7").Value = 4.974013
$ws.Range("O7").Value = 0.9482949895113968
$ws.Range("P7").Value = 0.9367771259737823
$ws.Range("Q7").Value = 7.954222733028001
$ws.Range("R7").Value = 47.725336398168
$ws.Range("S7").Value = 0.1583181012337233
$ws.Range("T7").Value = 0.1761701351958673
$ws.Range("G8").Value = 2.617047
$ws.Range("H8").Value = 7.851141
$ws.Range("I8").Value = 0.1366085279147753
$ws.Range("J8").Value = 0.1538816040404024
$ws.Range("M8").Value = 0.07111099999999999
$ws.Range("N8").Value = 0.142222
$ws.Range("O8").Value = 0.02711460746047303
$ws.Range("P8").Value = 0.02678527708115022
$ws.Range("Q8").Value = 0.186100829217
$ws.Range("R8").Value = 1.116604975302
$ws.Range("S8").Value = 0.003704086610162206
$ws.Range("T8").Value = 0.004121761401914024
$ws.Range("G9").Value = 2.617047
$ws.Range("H9").Value = 7.851141
$ws.Range("I9").Value = 0.1366085279147753
$ws.Range("J9").Value = 0.1538816040404024
$ws.Range("O9").Value = 0.0245904030281302
$ws.Range("P9").Value = 0.03643759694506741
$ws.Range("Q9").Value = 0.168775978077
$ws.Range("R9").Value = 1.518983802693
$ws.Range("S9").Value = 0.003359258758503901
$ws.Range("T9").Value = 0.00560707586528464
$ws.Range("G10").Value = 2.617047
$ws.Range("H10").Value = 7.851141
$ws.Range("I10").Value = 0.1366085279147753
$ws.Range("J10").Value = 0.1538816040404024
$ws.Range("M10").Value = 2.4870065
$ws.Range("N10").Value = 4.974013
$ws.Range("O10").Value = 0.9482949895113968
$ws.Range("P10").Value = 0.9367771259737823
$ws.Range("Q10").Value = 6.5086128998055
$ws.Range("R10").Value = 39.051677398833
$ws.Range("S10").Value = 0.1295451825461092
$ws.Range("T10").Value = 0.1441527667732038
$ws.Range("G11").Value = 4.784714
$ws.Range("H11").Value = 9.569428
$ws.Range("I11").Value = 0.2497596474320929
$ws.Range("J11").Value = 0.1875598630045162
$ws.Range("M11").Value = 0.07111099999999999
$ws.Range("N11").Value = 0.142222
$ws.Range("O11").Value = 0.02711460746047303
$ws.Range("P11").Value = 0.02678527708115022
$ws.Range("Q11").Value = 0.340245797254
$ws.Range("R11").Value = 1.360983189016
$ws.Range("S11").Value = 0.006772134799587341
$ws.Range("T11").Value = 0.005023842899878543
$ws.Range("G12").Value = 4.784714
$ws.Range("H12").Value = 9.569428
$ws.Range("I12").Value = 0.2497596474320929
$ws.Range("J12").Value = 0.1875598630045162
$ws.Range("O12").Value = 0.0245904030281302
$ws.Range("P12").Value = 0.03643759694506741
$ws.Range("Q12").Value = 0.3085709905740001
$ws.Range("R12").Value = 1.851425943444
$ws.Range("S12").Value = 0.006141690390518868
$ws.Range("T12").Value = 0.006834230691230621
$ws.Range("G13").Value = 4.784714
$ws.Range("H13").Value = 9.569428
$ws.Range("I13").Value = 0.2497596474320929
$ws.Range("J13").Value = 0.1875598630045162
$ws.Range("M13").Value = 2.4870065
$ws.Range("N13").Value = 4.974013
$ws.Range("O13").Value = 0.9482949895113968
$ws.Range("P13").Value = 0.9367771259737823
$ws.Range("Q13").Value = 11.899614818641
$ws.Range("R13").Value = 47.59845927456401
$ws.Range("S13").Value = 0.2368458222419867
$ws.Range("T13").Value = 0.175701789413407
$ws.Range("G14").Value = 4.899255333333334
$ws.Range("H14").Value = 14.697766
$ws.Range("I14").Value = 0.2557386470190557
$ws.Range("J14").Value = 0.2880747916628283
$ws.Range("M14").Value = 0.07111099999999999
$ws.Range("N14").Value = 0.142222
$ws.Range("O14").Value = 0.02711460746047303
$ws.Range("P14").Value = 0.02678527708115022
$ws.Range("Q14").Value = 0.3483909460086667
$ws.Range("R14").Value = 2.090345676052
$ws.Range("S14").Value = 0.006934253026394167
$ws.Range("T14").Value = 0.007716163114783479
$ws.Range("G15").Value = 4.899255333333334
$ws.Range("H15").Value = 14.697766
$ws.Range("I15").Value = 0.2557386470190557
$ws.Range("J15").Value = 0.2880747916628283
$ws.Range("O15").Value = 0.0245904030281302
$ws.Range("P15").Value = 0.03643759694506741
$ws.Range("Q15").Value = 0.3159578757020001
$ws.Range("R15").Value = 2.843620881318
$ws.Range("S15").Value = 0.006288716400067307
$ws.Range("T15").Value = 0.0104967531486444
$ws.Range("G16").Value = 4.899255333333334
$ws.Range("H16").Value = 14.697766
$ws.Range("I16").Value = 0.2557386470190557
$ws.Range("J16").Value = 0.2880747916628283
$ws.Range("M16").Value = 2.4870065
$ws.Range("N16").Value = 4.974013
$ws.Range("O16").Value = 0.9482949895113968
$ws.Range("P16").Value = 0.9367771259737823
$ws.Range("Q16").Value = 12.18447985915967
$ws.Range("R16").Value = 73.10687915495801
$ws.Range("S16").Value = 0.2425156775925942
$ws.Range("T16").Value = 0.2698618753994004
$ws.Range("G17").Value = 1.991497666666667
$ws.Range("H17").Value = 5.974493
$ws.Range("I17").Value = 0.1039551695437809
$ws.Range("J17").Value = 0.117099484796943
$ws.Range("M17").Value = 0.07111099999999999
$ws.Range("N17").Value = 0.142222
$ws.Range("O17").Value = 0.02711460746047303
$ws.Range("P17").Value = 0.02678527708115022
$ws.Range("Q17").Value = 0.1416173905743333
$ws.Range("R17").Value = 0.8497043434459999
$ws.Range("S17").Value = 0.002818703615666542
$ws.Range("T17").Value = 0.003136542146346057
$ws.Range("G18").Value = 1.991497666666667
$ws.Range("H18").Value = 5.974493
$ws.Range("I18").Value = 0.1039551695437809
$ws.Range("J18").Value = 0.117099484796943
$ws.Range("O18").Value = 0.0245904030281302
$ws.Range("P18").Value = 0.03643759694506741
$ws.Range("Q18").Value = 0.128433676021
$ws.Range("R18").Value = 1.155903084189
$ws.Range("S18").Value = 0.002556299515939179
$ws.Range("T18").Value = 0.004266823829506059
$ws.Range("G19").Value = 1.991497666666667
$ws.Range("H19").Value = 5.974493
$ws.Range("I19").Value = 0.1039551695437809
$ws.Range("J19").Value = 0.117099484796943
$ws.Range("M19").Value = 2.4870065
$ws.Range("N19").Value = 4.974013
$ws.Range("O19").Value = 0.9482949895113968
$ws.Range("P19").Value = 0.9367771259737823
$ws.Range("Q19").Value = 4.952867641734834
$ws.Range("R19").Value = 29.717205850409
$ws.Range("S19").Value = 0.09858016641217522
$ws.Range("T19").Value = 0.1096961188210909
